$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update table text (row/column labels renamed to reflect updated data
# releases: employment volumes/occupation split, FE achievements scope,
# enterprises-by-industry rename, and refreshed release period) ---
$ws.Range("A2").Value = "Employment volumes"
$ws.Range("A3").Value = "Employment by occupation"
$ws.Range("A5").Value = "Further education and skills achievements and participation by provision, level and age group"
$ws.Range("D5").Value = "Aug 2022 – Jul 2023 (Nov 23)"
$ws.Range("D6").Value = "Aug 2022 – Jul 2023 (Nov 23)"
$ws.Range("A9").Value = "Enterprises by employment industry"

# --- Tidy up view state: scroll/selection left how the author last saved it ---
$win = $excel.ActiveWindow
$win.ScrollColumn = 2
$win.ScrollRow = 1
$ws.Range("D7").Select()
